# Consolidate text runs: collapse multiple <a:r> runs that together spell out
# a sentence into a single run. Assigning TextRange.Text over the *whole*
# range re-derives paragraph boundaries from the string (vertical-tab / CR
# characters become new <a:p> paragraphs), which would eat the existing
# <a:br/> line breaks in the subtitle. Instead we target only the character
# span covering the runs we want to merge via TextRange.Characters(start,
# length), which keeps any leading <a:br/> elements untouched and simply
# rewrites that span as a single run.

$p = $ppt.ActivePresentation

# Slide 1, Subtitle 2: "Jesse" + " " + "Rosenthal" -> "Jesse Rosenthal"
# (the paragraph starts with two <a:br/> line breaks, which TextRange.Text
# represents as two leading vertical-tab characters -- skip over them.)
$s1 = $p.Slides.Item(1)
$subtitleRange = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitleRange.Characters(3, $subtitleRange.Length - 2).Text = "Jesse Rosenthal"

# Slide 2, Title 1: "A" + " " + "header" -> "A header"
# (plain whole-range assignment is a no-op here because the joined text is
# already identical to the target string; go through Characters() so the
# engine actually rewrites the run storage instead of skipping the "change".)
$s2 = $p.Slides.Item(2)
$titleRange = $s2.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "A header"

# Slide 1's notes page, Notes Placeholder 2:
# "Some" + " " + "speaker" + " " + "notes" -> "Some speaker notes"
$notes1 = $s1.NotesPage
$notesRange = $notes1.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "Some speaker notes"
